$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 458; this pushes the existing
# rows 458-471 down to 461-474 (and copies the row-458 formatting,
# notably the date-number-format style on column D, into the new rows).
$ws.Rows("458:460").Insert()

# Common (constant) values shared by every row in this data block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$cultivar  = "Cultivar IV Región"
$unidad    = '$/bandeja 10 kilos'
$provinciaElqui = "Provincia del Elquí"
$numObs    = 10

function Set-ChirimoyaRow($row, $fecha, $calidad, $nroObs, $precioMin, $precioMax, $precioMod, $provincia, $precioKilo) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $cultivar
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $nroObs
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioMod
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $provincia
    $ws.Cells.Item($row, 19).Value = $precioKilo
    $ws.Cells.Item($row, 20).Value = $numObs
}

Set-ChirimoyaRow 458 45239 "Especial" 45 28000 28000 28000 $provinciaElqui 2800
Set-ChirimoyaRow 459 45239 "Primera"  50 25000 25000 25000 $provinciaElqui 2500
Set-ChirimoyaRow 460 45239 "Segunda"  40 20000 20000 20000 $provinciaElqui 2000
